# "update scripts wuth new tpm"
#
# The NATMI pipeline was re-run with new TPM-derived expression numbers.
# The "ECs" sending-cluster (old rows 2 & 3) dropped out of this
# particular sheet entirely, and the remaining "FAPs" / "Resolving-Mac"
# sending-cluster rows (old rows 4-7) moved up to rows 2-5 and got
# refreshed edge-weight / specificity figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two "ECs" sending-cluster rows; everything below shifts up.
$ws.Rows("2:3").Delete()

# Row 2: FAPs -> Siglec1/Spn -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Siglec1"
$ws.Range("C2").Value = "Spn"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03801766666666666
$ws.Range("H2").Value = 0.114053
$ws.Range("I2").Value = 0.001128449675396954
$ws.Range("J2").Value = 0.001128449675396954
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1149353333333333
$ws.Range("N2").Value = 0.344806
$ws.Range("O2").Value = 0.05719122335670149
$ws.Range("P2").Value = 0.05719122335670149
$ws.Range("Q2").Value = 0.004369573190888889
$ws.Range("R2").Value = 0.039326158718
$ws.Range("S2").Value = 0.00006453741743242446
$ws.Range("T2").Value = 0.00006453741743242446

# Row 3: FAPs -> Siglec1/Spn -> Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Siglec1"
$ws.Range("C3").Value = "Spn"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03801766666666666
$ws.Range("H3").Value = 0.114053
$ws.Range("I3").Value = 0.001128449675396954
$ws.Range("J3").Value = 0.001128449675396954
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.894732
$ws.Range("N3").Value = 5.684196
$ws.Range("O3").Value = 0.9428087766432985
$ws.Range("P3").Value = 0.9428087766432984
$ws.Range("Q3").Value = 0.07203328959866667
$ws.Range("R3").Value = 0.648299606388
$ws.Range("S3").Value = 0.001063912257964529
$ws.Range("T3").Value = 0.001063912257964529

# Row 4: Resolving-Mac -> Siglec1/Spn -> ECs
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Siglec1"
$ws.Range("C4").Value = "Spn"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 33.652157
$ws.Range("H4").Value = 100.956471
$ws.Range("I4").Value = 0.998871550324603
$ws.Range("J4").Value = 0.9988715503246031
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1149353333333333
$ws.Range("N4").Value = 0.344806
$ws.Range("O4").Value = 0.05719122335670149
$ws.Range("P4").Value = 0.05719122335670149
$ws.Range("Q4").Value = 3.867821882180666
$ws.Range("R4").Value = 34.810396939626
$ws.Range("S4").Value = 0.05712668593926906
$ws.Range("T4").Value = 0.05712668593926906

# Row 5: Resolving-Mac -> Siglec1/Spn -> Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Siglec1"
$ws.Range("C5").Value = "Spn"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 33.652157
$ws.Range("H5").Value = 100.956471
$ws.Range("I5").Value = 0.998871550324603
$ws.Range("J5").Value = 0.9988715503246031
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.894732
$ws.Range("N5").Value = 5.684196
$ws.Range("O5").Value = 0.9428087766432985
$ws.Range("P5").Value = 0.9428087766432984
$ws.Range("Q5").Value = 63.76181873692399
$ws.Range("R5").Value = 573.8563686323159
$ws.Range("S5").Value = 0.9417448643853339
$ws.Range("T5").Value = 0.9417448643853339
